# Apply the "Add data for 2022-02-21" update:
#  - rename the sheet / tab from "Through 2022-02-12" to "Through 2022-02-13"
#  - update the header cell (I1) text from "2022 (through 02-12)" to "2022 (through 02-13)"
#  - bump I3 (February total) from 56 to 59
#  - bump I14 (Total row) from 217 to 220

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (this updates both the sheet tab name and workbook.xml <sheet name=.../>)
$ws.Name = "Through 2022-02-13"

# Update the "2022 (through ...)" header label in column I, row 1
$ws.Range("I1").Value = "2022 (through 02-13)"

# Update the affected numeric cells
$ws.Cells.Item(3, 9).Value = 59
$ws.Cells.Item(14, 9).Value = 220
